$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 08.02.2022 09:45"

# 2. D6: change from text "+0.6" to numeric value 0.6
$ws.Range("D6").Value = 0.6

# 3. E6: change from text "2022-02-08 09:30:21" to the equivalent numeric date
#    serial value, matching the style/number format used by the other date
#    cells in column E (E2:E10, style s="2")
$ws.Range("E6").Value = 44600.39607638889
$ws.Range("E6").NumberFormat = $ws.Range("E7").NumberFormat
